$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: update Objetivos (PT) text in B and C
$ws.Range("B10").Value = "Fornecer fundamentos teóricos de limite e derivadas, destacando aspectos geométricos e interpretações físicas, elementos fundamentais para estudos de Engenharia"
$ws.Range("C10").Value = "Fornecer fundamentos teóricos de limite e derivadas, destacando aspectos geométricos e interpretações físicas, elementos fundamentais para estudos de Engenharia"

# Row 13: remove label (col A), update B/C to teacher name, clear custom row height
$ws.Range("A13").ClearContents()
$ws.Range("B13").Value = "5840692 - Diovana Aparecida dos Santos Napoleão"
$ws.Range("C13").Value = "5840692 - Diovana Aparecida dos Santos Napoleão"
$ws.Rows.Item(13).AutoFit()

# Row 14: Programa resumido (label+text)
$ws.Range("A14").Value = "Programa resumido:"
$ws.Range("B14").Value = "Números Reais, funções de variável real, limites e derivadas de funções Reais. Aplicações da derivada e Fórmula de Taylor."
$ws.Range("C14").Value = "Números Reais, funções de variável real, limites e derivadas de funções Reais. Aplicações da derivada e Fórmula de Taylor."

# Row 15: Short syllabus (label+text), row height 120 -> 60
$ws.Range("A15").Value = "Short syllabus:"
$ws.Range("B15").Value = "Descrição do programa resumido em inglês.Real numbers, real functions, limits and derivatives of real functions.  Applications of the derivative and Taylor’s Formula."
$ws.Range("C15").Value = "Descrição do programa resumido em inglês.Real numbers, real functions, limits and derivatives of real functions.  Applications of the derivative and Taylor’s Formula."
$ws.Rows.Item(15).RowHeight = 60

# Row 16: Programa (label+text)
$ws.Range("A16").Value = "Programa:"
$ws.Range("B16").Value = "•Números e Funções Reais: função trigonométrica, exponencial e logarítmica. Função composta e inversa.•Limite: Definição, propriedades algébricas e Teorema do confronto. Limites infinitos e ao infinito.•Continuidade de funções Reais: Teorema de Weierstrass e teorema do valor intermediário.•Derivada de funções Reais: Definição, Interpretação física e geométrica, regras de derivação, regra da cadeia, derivada da função inversa e derivação implícita, Regra de l’ hopital, Teorema do valor Médio e consequências, Formula de Taylor, taxas de variação, máximos e mínimos (otimização)."
$ws.Range("C16").Value = "•Números e Funções Reais: função trigonométrica, exponencial e logarítmica. Função composta e inversa.•Limite: Definição, propriedades algébricas e Teorema do confronto. Limites infinitos e ao infinito.•Continuidade de funções Reais: Teorema de Weierstrass e teorema do valor intermediário.•Derivada de funções Reais: Definição, Interpretação física e geométrica, regras de derivação, regra da cadeia, derivada da função inversa e derivação implícita, Regra de l’ hopital, Teorema do valor Médio e consequências, Formula de Taylor, taxas de variação, máximos e mínimos (otimização)."

# Row 17: Syllabus (label updated), add new B/C text, set row height to 120
$ws.Range("A17").Value = "Syllabus:"
$ws.Range("B16").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("B17").Value = "•Real Numbers and Real Functions: trigonometric, exponential and logarithmic functions. Composite and inverse functions.•Limits: Definition, algebraic properties and squeeze theorem. Infinite limits and Limits to infinite.•Continuity: Weierstrass theorem and intermediate value theorem.•Derivative of real functions: Definition, geometrical and physics interpretations, derivative rules, chain rule, derivative of inverse and implicit functions, l’hopital rule, mean value theorem and consequences, Taylor’s Formula,  Maximum and Minimum Problems"
$ws.Range("C16").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C17").Value = "•Real Numbers and Real Functions: trigonometric, exponential and logarithmic functions. Composite and inverse functions.•Limits: Definition, algebraic properties and squeeze theorem. Infinite limits and Limits to infinite.•Continuity: Weierstrass theorem and intermediate value theorem.•Derivative of real functions: Definition, geometrical and physics interpretations, derivative rules, chain rule, derivative of inverse and implicit functions, l’hopital rule, mean value theorem and consequences, Taylor’s Formula,  Maximum and Minimum Problems"
$ws.Rows.Item(17).RowHeight = 120

# Row 18: Avaliacao (label only) - clear old B/C content, remove custom row height
$ws.Range("A18").Value = "Avaliação:"
$ws.Range("B18").ClearContents()
$ws.Range("C18").ClearContents()
$ws.Rows.Item(18).AutoFit()

# Row 19: Metodo (label updated)
$ws.Range("A19").Value = "Método:"

# Row 20: Criterio (label updated)
$ws.Range("A20").Value = "Critério:"

# Row 21: Norma de recuperacao (label updated), row height 120 -> 60
$ws.Range("A21").Value = "Norma de recuperação:"
$ws.Rows.Item(21).RowHeight = 60

# Row 22: new Bibliografia row - copy formatting from row 21, then set values, row height 120
$ws.Range("A21").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("A22").Value = "Bibliografia:"
$ws.Range("B21").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("B22").Value = "STEWART, James. Cálculo São Paulo: Cengage Learning, 2009. v.1.`nANTON, Howard. Cálculo: um novo horizonte. Porto Alegre: Bookman, 2007.`nTHOMAS, George B. Cálculo São Paulo: Pearson Addison  Wesley, 2009. v.1,`nGUIDORIZZI, Hamilton. Um curso de cálculo. Rio de Janeiro: Livros Técnicos e Científicos, 2001. v.1.`nFLEMMING, Diva M.; GONÇALVES, Mirian B. Cálculo A. São Paulo: Pearson Prentice Hall, 2009."
$ws.Range("C21").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = "STEWART, James. Cálculo São Paulo: Cengage Learning, 2009. v.1.`nANTON, Howard. Cálculo: um novo horizonte. Porto Alegre: Bookman, 2007.`nTHOMAS, George B. Cálculo São Paulo: Pearson Addison  Wesley, 2009. v.1,`nGUIDORIZZI, Hamilton. Um curso de cálculo. Rio de Janeiro: Livros Técnicos e Científicos, 2001. v.1.`nFLEMMING, Diva M.; GONÇALVES, Mirian B. Cálculo A. São Paulo: Pearson Prentice Hall, 2009."
$ws.Rows.Item(22).RowHeight = 120

